$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new column before C (shifts old C..I -> D..J) ---
$ws.Columns.Item(3).Insert()

# --- 2. Header row (row 1): new C1 header, rest already correct text after shift ---
$ws.Range("C1").Value = "Thuong hieu"
$ws.Range("C1").Value = "Thương hiệu"

# --- 3. Row 2: clear stale cells, then set fresh values ---
$ws.Range("F2").ClearContents()
$ws.Range("H2").ClearContents()

$ws.Range("A2").Value = "B05-10.6-Standard"
$ws.Range("B2").Value = "Standard"
$ws.Range("C2").Value = "Fargo"

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "10.6"
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").Value = "B"

$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "05"

$ws.Range("I2").Value = "Cai"
$ws.Range("I2").Value = "Cái"
$ws.Range("J2").Value = "khong"
$ws.Range("J2").Value = "không"

# --- 4. Row 3: replace entirely ---
$ws.Range("A3").Value = "B06-10.6-Standard Isee B&L"
$ws.Range("B3").Value = "Standard Isee B&L"
$ws.Range("C3").Value = "Isee B&L"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "10.6"
$ws.Range("D3").Style = "Normal"

$ws.Range("E3").Value = "B"

$ws.Range("G1").NumberFormat = "@"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "06"

$ws.Range("I3").Value = "Cái"
$ws.Range("J3").Value = "không"

# --- 5. Row 4: brand new row ---
$ws.Range("A4").Value = "Dung dịch thử nghiệm"
$ws.Range("B4").Value = "Dung dịch"
$ws.Range("C4").Value = "Fargo"
$ws.Range("I4").Value = "Chai"
$ws.Range("J4").Value = "không"

# --- 6. Column widths (best effort given engine rounding to 1/6 character units) ---
$ws.Columns.Item(1).ColumnWidth = 28.75
$ws.Columns.Item(2).ColumnWidth = 15.59
$ws.Columns.Item(3).ColumnWidth = 13.42
$ws.Columns.Item(4).ColumnWidth = 11.09
$ws.Columns.Item(5).ColumnWidth = 9.09
$ws.Columns.Item(6).ColumnWidth = 6.09
$ws.Columns.Item(7).ColumnWidth = 8.25
$ws.Columns.Item(10).ColumnWidth = 10.59

# --- 7. Selection ---
$ws.Range("D7").Select()

Write-Output "done"
